$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.788.35"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.828.82"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'230.79"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'39.26"
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").Value = "'0.327"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "'0.0989"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").Value = "2.095.02"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "'11.31"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "1.828.13"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "'0.669"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").Value = "34.769.45"
$ws.Range("D18").Value = "'69.47"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "0.0₃0785"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "'239.96"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "'12.15"
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'171.87"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("D26").Value = "'7.75"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("D27").Value = "'0.124"
$ws.Range("E27").Value = "  +2.52%  "
$ws.Range("D28").Value = "'17.32"
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("E29").Value = "  -8.09%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "'0.0550"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "'3.91"
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("E33").Value = "  -0.96%  "
$ws.Range("D34").Value = "'1.85"
$ws.Range("E34").Value = "  +3.74%  "
$ws.Range("E35").Value = "  +7.32%  "
$ws.Range("D36").Value = "'1.42"
$ws.Range("E36").Value = "  +11.09%  "
$ws.Range("D37").Value = "'0.698"
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("D38").Value = "'91.23"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").Value = "'1.04"
$ws.Range("E39").Value = "  +5.72%  "
$ws.Range("D40").Value = "1.341.64"
$ws.Range("E40").Value = "  +2.81%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'14.48"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("E44").Value = "  -3.95%  "
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").Value = "'6.25"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").Value = "2.010.44"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  +4.09%  "
$ws.Range("D51").Value = "'98.44"
$ws.Range("E51").Value = "  -0.31%  "
